$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: remove the "Our model uses the Euler method ..." sentence from
# the first paragraph, keeping the final ". " after "...reasonable amount of
# time".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " Our model uses the Euler method to solve the differential equations. This could be replaced by a more accurate method like Runge Kutta but we considered Euler as precise enough.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Re-split the remaining run into "...can be done" | " in a reasonable
# amount of time" | ". " the way the authored commit shows it, by nudging a
# direct formatting property on a sub-range (forces the run boundary without
# altering the visible formatting).
$r1 = $d.Content
$r1.Find.Execute(" in a reasonable amount of time")
$s1 = $r1.Start
$e1 = $r1.End
$splitA = $d.Range($s1, $e1)
$splitA.Font.Bold = $true
$splitA.Font.Bold = $false

# ---------------------------------------------------------------------------
# Change 2: split "...certainly occurs in nature could" into three runs
# around the word "occurs" (mirrors the proofing-mark wrap in the commit).
# Use the unique surrounding phrase so we land on this "occurs", not the
# other one earlier in the document.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("certainly occurs in nature could")
$phraseStart = $r2.Start
$occursStart = $phraseStart + 10
$occursEnd = $occursStart + 6

$split2A = $d.Range($occursStart, $occursEnd)
$split2A.Font.Bold = $true
$split2A.Font.Bold = $false

# ---------------------------------------------------------------------------
# Change 3: "It is hard to analyse ... " paragraph: wrap the first "preys"
# in its own run, and wrap the "a" in "is a interplay" in its own run, the
# way the authored commit shows it.
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("behaviour of preys")
$phrase3Start = $r3.Start
$preysStart = $phrase3Start + 13
$preysEnd = $preysStart + 5
$split3A = $d.Range($preysStart, $preysEnd)
$split3A.Font.Bold = $true
$split3A.Font.Bold = $false

$r4 = $d.Content
$r4.Find.Execute("Predator prey behaviour is a interplay")
$phrase4Start = $r4.Start
$aStart = $phrase4Start + 27
$aEnd = $aStart + 1
$split4A = $d.Range($aStart, $aEnd)
$split4A.Font.Bold = $true
$split4A.Font.Bold = $false
